# Add the bash entry on gpg: append a new row (62) to the first worksheet
# with the topic "gpg/ pgp" and widen column C to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$gpgBody = "# What is pgp/ gpg`n" + `
  "It is an util that implement RSA encryption and decryption methodology`n" + `
  "# Installation`n" + `
  "apt-cyg/ apt-get install gnupg `n" + `
  "# Operation flow - Create my key-pack and export my public key `n" + `
  "    * gpg --gen-key  //gpg will ask you a name and it will be regarded as the key id in the local key db`n" + `
  "    * gpg --list-key`n" + `
  "    * gpg --output ""myles_public_key.asc"" --export ""mykey"" //""mykey"" is the key id in local key databse`n" + `
  "# Operation flow - Decrypt file that is encrypted by others with my public key`n" + `
  "    * gpg --decrypt-files ""file_encrypt_with_my_public_key"" //gpg will find a key in db that can decrypt`n" + `
  "# Operation flow - Import other's public key and use it to encrypt file `n" + `
  "    * gpg --import {public_key_from_other.sac}`n" + `
  "    * gpg --list-key  //In here learnt the key id of the newly imported key`n" + `
  "    * gpg --recipient ""their_keyid"" --output ""outputfilename.gpg"" --encrypt ""file_to_be_encrypted"""

$ws.Range("A62").Value = "gpg/ pgp"
$ws.Range("B62").Value = "# Basic Opts"
$ws.Range("C62").Value = $gpgBody

# Match the surrounding rows' formatting: wrapped text in column C.
$ws.Range("C62").WrapText = $true
$ws.Rows.Item(62).RowHeight = 220.5

# Column C widens to fit the long new entry.
$ws.Columns.Item(3).ColumnWidth = 86.833333

# Leave the sheet scrolled/selected at the newly added row, like the author did.
$ws.Range("C63").Select() | Out-Null
